# Daily attendance processing - 2025-12-13 17:26:33
# Normalizes the "Recorded By" (column G) entries so that the "System"
# (and lowercase "system") automated-recorder tokens are listed before
# any real/backup email addresses, matching the canonical ordering used
# by the attendance pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact before -> after mappings observed for the "Recorded By" column.
$map = @{
    "backup@backdoor.com, system, System" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "System, admin@admin.com"             = "admin@admin.com, System"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value()
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
